$d = $word.ActiveDocument

$replacements = @(
    @("43×91=", "78×82="),
    @("88×66=", "78×51="),
    @("67×27=", "71×22="),
    @("47×48=", "29×18="),
    @("91×91=", "79×27="),
    @("28×49=", "99×99="),
    @("52×70=", "63×16="),
    @("90×95=", "65×28="),
    @("88×77=", "18×67="),
    @("99×48=", "81×17="),
    @("21×62=", "50×48="),
    @("88×54=", "34×34="),
    @("63×38=", "74×15="),
    @("30×41=", "33×56="),
    @("23×61=", "75×97="),
    @("28×37=", "78×30="),
    @("34×41=", "76×99="),
    @("51×25=", "36×60="),
    @("31×57=", "32×53="),
    @("75×91=", "63×43="),
    @("73×40=", "50×78="),
    @("63×74=", "90×22="),
    @("29×76=", "12×90="),
    @("62×71=", "95×96="),
    @("57×62=", "95×36=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
